$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos value becomes the professor's name line ---
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# --- Row 13: gains label "Programa resumido:" in col A, value becomes "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: label becomes "Short syllabus:", value becomes the short-syllabus English text ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Management and strategic process; Strategic, tactical and operational planning; Organizational guidelines: mission, vision and objectives; Strategic formulation; Strategy implementation; Strategic control; Planning of organizational units; Strategic management applied."
$ws.Range("C14").Value = "Management and strategic process; Strategic, tactical and operational planning; Organizational guidelines: mission, vision and objectives; Strategic formulation; Strategy implementation; Strategic control; Planning of organizational units; Strategic management applied."

# --- Row 15: label becomes "Programa:", value becomes "01/01/2021", height 120 ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"
$ws.Rows(15).RowHeight = 120

# --- Row 16: label becomes "Syllabus:", value becomes the English syllabus text (content unchanged, height stays 120) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Motivations and Challenges for the strategy; 2. Basic concepts of strategy; 3. Strategic Management; 4. Strategic Transformation; 5. Analysis of the External Environment; 6. Analysis of Turbulence and Vulnerability; 7. Analysis of the Internal Environment; 8. Representation of the Portfolio; 9. Portfolio Balancing Strategy; 10. Formulation of Strategies; 11. Strategic Training; 12. The Strategic Plan; 13. Strategic Planning Methodology; 14. Strategic Planning Workshop; 15. Implementation of Strategic Management."
$ws.Range("C16").Value = "1. Motivations and Challenges for the strategy; 2. Basic concepts of strategy; 3. Strategic Management; 4. Strategic Transformation; 5. Analysis of the External Environment; 6. Analysis of Turbulence and Vulnerability; 7. Analysis of the Internal Environment; 8. Representation of the Portfolio; 9. Portfolio Balancing Strategy; 10. Formulation of Strategies; 11. Strategic Training; 12. The Strategic Plan; 13. Strategic Planning Methodology; 14. Strategic Planning Workshop; 15. Implementation of Strategic Management."

# --- Row 17: label becomes "Avaliação:"; B/C cleared (row no longer carries a value); height reverts to default ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").ClearContents()
$ws.Rows(17).RowHeight = 15

# --- Row 18: label becomes "Método:", value becomes the professor's name line, height 60 ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows(18).RowHeight = 60

# --- Row 19: label becomes "Critério:" (value text unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label becomes "Norma de recuperação:" (value text unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label becomes "Bibliografia:" (value text unchanged), height becomes 120 ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# --- Row 22: label becomes "Requisitos:"; B/C cleared; height reverts to default ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").ClearContents()
$ws.Rows(22).RowHeight = 15

# --- Row 23: A cleared (no longer has a label); value becomes the LOQ4239 requirement line, height 30 ---
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOQ4239 -  Administração e Organização I  (Requisito fraco)" + [char]10
$ws.Range("C23").Value = "LOQ4239 -  Administração e Organização I  (Requisito fraco)" + [char]10
$ws.Rows(23).RowHeight = 30

# --- Row 24: value becomes the LOQ4240 requirement line (shifted up from row 25) ---
$ws.Range("B24").Value = "LOQ4240 -  Administração e Organização II  (Requisito fraco)" + [char]10
$ws.Range("C24").Value = "LOQ4240 -  Administração e Organização II  (Requisito fraco)" + [char]10

# --- Row 25 no longer exists: delete it outright (shifts nothing below, fixes the used range to C24) ---
$ws.Rows(25).Delete()
